# Automatic update of files.
# Update the "Förändrad" (changed/updated) date column (C) for rows 2-10
# from 2023-09-16 (45185) to 2023-10-05 (45204).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C10").Value = 45204
